$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.156.69'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.372.58'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.79'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.32'
$ws.Range('E6').Value = '  +1.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.503'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.483'
$ws.Range('E9').Value = '  -2.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.47'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('E11').Value = '  +3.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0790'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.29'
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.83'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.736.83'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.409.90'
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.805'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.170.60'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.99'
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.31'
$ws.Range('E20').Value = '  +1.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.98'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.72'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.48'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('E29').Value = '  +1.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.95'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('B33').Value = 'Celestia'
$ws.Range('C33').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.86'
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.111'
$ws.Range('E34').Value = '  +10.69%  '
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '127.49'
$ws.Range('E36').Value = '  +14.31%  '
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.86'
$ws.Range('E38').Value = '  +3.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.34'
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.26'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.04'
$ws.Range('E42').Value = '  -6.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.930.98'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.22'
$ws.Range('E47').Value = '  -8.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.594.60'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.53'
$ws.Range('E49').Value = '  +3.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.70'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.77'
$ws.Range('E51').Value = '  -2.67%  '
